# cryptos.xlsx price/volume refresh -- commit: "Updated cryptos list on Mon Jun  5 21:40:09 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells in this sheet are stored as text (t="inlineStr"),
# including values that look numeric (e.g. "1.002", "19.82") or date-like
# (e.g. "25.583.91"). Prefixing with a literal leading apostrophe forces
# Excel/COM to keep them as text instead of auto-converting to numbers/dates.

$ws.Range("D2").Value = "'" + '25.583.91'
$ws.Range("E2").Value = "'" + '  -6.03%  '
$ws.Range("D3").Value = "'" + '1.801.19'
$ws.Range("E3").Value = "'" + '  -5.45%  '
$ws.Range("D4").Value = "'" + '1.002'
$ws.Range("E4").Value = "'" + '  +0.06%  '
$ws.Range("D5").Value = "'" + '273.47'
$ws.Range("E5").Value = "'" + '  -10.68%  '
$ws.Range("D7").Value = "'" + '0.4998'
$ws.Range("E7").Value = "'" + '  -7.78%  '
$ws.Range("D8").Value = "'" + '0.3491'
$ws.Range("E8").Value = "'" + '  -8.30%  '
$ws.Range("D9").Value = "'" + '43.74'
$ws.Range("E9").Value = "'" + '  -4.96%  '
$ws.Range("D10").Value = "'" + '0.06575'
$ws.Range("E10").Value = "'" + '  -9.80%  '
$ws.Range("D11").Value = "'" + '19.82'
$ws.Range("D12").Value = "'" + '0.8310'
$ws.Range("E12").Value = "'" + '  -7.90%  '
$ws.Range("D13").Value = "'" + '0.07768'
$ws.Range("E13").Value = "'" + '  -5.27%  '
$ws.Range("D14").Value = "'" + '1.799.29'
$ws.Range("E14").Value = "'" + '  +44.92%  '
$ws.Range("D15").Value = "'" + '5.034'
$ws.Range("E15").Value = "'" + '  -5.90%  '
$ws.Range("D16").Value = "'" + '87.03'
$ws.Range("E16").Value = "'" + '  -9.02%  '
$ws.Range("D17").Value = "'" + '1.002'
$ws.Range("E17").Value = "'" + '  +0.03%  '
$ws.Range("D18").Value = "'" + '13.84'
$ws.Range("E18").Value = "'" + '  -6.65%  '
$ws.Range("E19").Value = "'" + '  +0.09%  '
$ws.Range("D20").Value = "'" + '0.000007922'
$ws.Range("E20").Value = "'" + '  -8.31%  '
$ws.Range("D21").Value = "'" + '25.664.10'
$ws.Range("D22").Value = "'" + '4.697'
$ws.Range("E22").Value = "'" + '  -6.95%  '
$ws.Range("D23").Value = "'" + '9.965'
$ws.Range("E23").Value = "'" + '  -7.82%  '
$ws.Range("D24").Value = "'" + '6.028'
$ws.Range("E24").Value = "'" + '  -7.44%  '
$ws.Range("D25").Value = "'" + '141.61'
$ws.Range("E25").Value = "'" + '  -4.52%  '
$ws.Range("D26").Value = "'" + '2.090'
$ws.Range("E26").Value = "'" + '  -9.09%  '
$ws.Range("E27").Value = "'" + '  -5.99%  '
$ws.Range("D28").Value = "'" + '16.86'
$ws.Range("E28").Value = "'" + '  -8.17%  '
$ws.Range("D29").Value = "'" + '107.74'
$ws.Range("E29").Value = "'" + '  -7.76%  '
$ws.Range("D30").Value = "'" + '4.296'
$ws.Range("E30").Value = "'" + '  -11.44%  '
$ws.Range("D31").Value = "'" + '4.170'
$ws.Range("E31").Value = "'" + '  -10.52%  '
$ws.Range("D32").Value = "'" + '0.08754'
$ws.Range("E32").Value = "'" + '  -4.91%  '
$ws.Range("E33").Value = "'" + '  -5.75%  '
$ws.Range("D34").Value = "'" + '2.882'
$ws.Range("E34").Value = "'" + '  -4.30%  '
$ws.Range("B35").Value = "'" + 'ARBITRUM'
$ws.Range("C35").Value = "'" + 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'" + '1.124'
$ws.Range("E35").Value = "'" + '  -7.87%  '
$ws.Range("B36").Value = "'" + 'ImmutableX'
$ws.Range("C36").Value = "'" + 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'" + '0.7158'
$ws.Range("E36").Value = "'" + '  -13.32%  '
$ws.Range("D37").Value = "'" + '1.000'
$ws.Range("E37").Value = "'" + '  +0.07%  '
$ws.Range("D38").Value = "'" + '3.019'
$ws.Range("E38").Value = "'" + '  -9.16%  '
$ws.Range("D39").Value = "'" + '0.01851'
$ws.Range("E39").Value = "'" + '  -7.22%  '
$ws.Range("D40").Value = "'" + '0.5103'
$ws.Range("E40").Value = "'" + '  -15.30%  '
$ws.Range("D41").Value = "'" + '2.252'
$ws.Range("E41").Value = "'" + '  -16.44%  '
$ws.Range("D42").Value = "'" + '0.9495'
$ws.Range("E42").Value = "'" + '  -11.66%  '
$ws.Range("D43").Value = "'" + '113.64'
$ws.Range("E43").Value = "'" + '  -1.98%  '
$ws.Range("D44").Value = "'" + '6.126'
$ws.Range("E44").Value = "'" + '  -7.97%  '
$ws.Range("D45").Value = "'" + '7.930'
$ws.Range("E45").Value = "'" + '  -14.30%  '
$ws.Range("E46").Value = "'" + '  +0.09%  '
$ws.Range("D47").Value = "'" + '0.1373'
$ws.Range("E47").Value = "'" + '  -10.28%  '
$ws.Range("D48").Value = "'" + '0.4517'
$ws.Range("E48").Value = "'" + '  -12.12%  '
$ws.Range("D49").Value = "'" + '9.246'
$ws.Range("E49").Value = "'" + '  -9.41%  '
$ws.Range("D50").Value = "'" + '35.81'
$ws.Range("E50").Value = "'" + '  -5.98%  '
$ws.Range("D51").Value = "'" + '1.479'
$ws.Range("E51").Value = "'" + '  -9.79%  '
